{"js": "// Redefinition of the unidirectional load transfer function for PI4 and PI5.\n// 1) Update the PI4 description text (appears both in the detailed PI section\n//    and in the condensed summary section; both copies share the same\n//    underlying sentence, so a single search/replace handles both).\nconst oldDescription = \"Time needed for unidirectional load transfer \\u2013 this PI is a scalar indicating the average time (across sts cycles) elapsed between the beginning of a sit to stand (or stand to sit) movement and the full transfer of the weight on a single platform. During sit-to-stand this time corresponds to the time needed to move the CoP purely on the ground platform. Data from both the Chair and lower limb kinematics are needed for calculating this PI. \";\nconst newDescription = \"Time needed for unidirectional load transfer \\u2013 this PI is an two elements array of scalars indicating the AP and ML unidirectional load transfer overshoot times, corresponding to the time at which the distance between the quiet standing CoP position and the local maxima of anteroposterior and medio-lateral CoP during sts transition are reached. The data is averaged across 5 STS cycles. Data from both the Chair and lower limb kinematics are needed for calculating this PI.\";\n\nconst descriptionMatches = context.document.body.search(oldDescription, { matchCase: true });\ndescriptionMatches.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < descriptionMatches.items.length; i++) {\n  descriptionMatches.items[i].insertText(newDescription, Word.InsertLocation.replace);\n}\nawait context.sync();\nif (descriptionMatches.items.length === 0) {\n  throw new Error(\"PI4 description text to replace was not found.\");\n}\n\n// 2) Rename the two corresponding-function file names (each name appears\n//    once in the detailed PI section and once in the summary section) to\n//    the new, shared function name \"unidirectional_load_transfer.m\".\nconst newFunctionName = \"unidirectional_load_transfer.m\";\n\nconst timeNeededMatches = context.document.body.search(\"time_needed_ult.m\", { matchCase: true });\ntimeNeededMatches.load(\"text\");\nconst ultOvershootMatches = context.document.body.search(\"ult_overshoot.m\", { matchCase: true });\nultOvershootMatches.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < timeNeededMatches.items.length; i++) {\n  timeNeededMatches.items[i].insertText(newFunctionName, Word.InsertLocation.replace);\n}\nfor (let i = 0; i < ultOvershootMatches.items.length; i++) {\n  ultOvershootMatches.items[i].insertText(newFunctionName, Word.InsertLocation.replace);\n}\nawait context.sync();\nif (timeNeededMatches.items.length === 0 || ultOvershootMatches.items.length === 0) {\n  throw new Error(\"One of the corresponding-function names to rename was not found.\");\n}\n", "ps1": "# Redefinition of the unidirectional load transfer function for PI4 and PI5.\n#\n# 1) Update the PI4 description text. The same underlying sentence appears\n#    twice in the document (once in the detailed PI4 section, once in the\n#    condensed summary section further down), so a single Find/Replace over\n#    the whole document body handles both occurrences.\n# 2) Rename the two \"Corresponding function\" file names (each of which also\n#    appears twice, in the detailed section and in the summary section) to\n#    the new, shared function name \"unidirectional_load_transfer.m\".\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $found) {\n        throw \"Text to replace was not found: $oldText\"\n    }\n}\n\n$oldDescription = \"Time needed for unidirectional load transfer \u2013 this PI is a scalar indicating the average time (across sts cycles) elapsed between the beginning of a sit to stand (or stand to sit) movement and the full transfer of the weight on a single platform. During sit-to-stand this time corresponds to the time needed to move the CoP purely on the ground platform. Data from both the Chair and lower limb kinematics are needed for calculating this PI. \"\n$newDescription = \"Time needed for unidirectional load transfer \u2013 this PI is an two elements array of scalars indicating the AP and ML unidirectional load transfer overshoot times, corresponding to the time at which the distance between the quiet standing CoP position and the local maxima of anteroposterior and medio-lateral CoP during sts transition are reached. The data is averaged across 5 STS cycles. Data from both the Chair and lower limb kinematics are needed for calculating this PI.\"\n\nReplace-AllText $oldDescription $newDescription\n\nReplace-AllText \"time_needed_ult.m\" \"unidirectional_load_transfer.m\"\nReplace-AllText \"ult_overshoot.m\" \"unidirectional_load_transfer.m\"\n"}
